$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cache the price/volume columns for the changed rows as text so Excel
# does not auto-convert values like "1.002" or "25.679.90" into numbers.
$priceVolRange = $ws.Range('D2:E51')
$priceVolRange.NumberFormat = '@'

$ws.Range('D2').Value = '25.679.90'
$ws.Range('E2').Value = '  -3.91%  '
$ws.Range('D3').Value = '1.743.73'
$ws.Range('E3').Value = '  -5.85%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '236.14'
$ws.Range('E5').Value = '  -10.40%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '0.4900'
$ws.Range('E7').Value = '  -9.14%  '
$ws.Range('D8').Value = '41.59'
$ws.Range('E8').Value = '  -7.81%  '
$ws.Range('D9').Value = '0.2494'
$ws.Range('E9').Value = '  -22.29%  '
$ws.Range('D10').Value = '0.05933'
$ws.Range('E10').Value = '  -16.23%  '
$ws.Range('D11').Value = '1.744.77'
$ws.Range('E11').Value = '  -5.85%  '
$ws.Range('D12').Value = '0.06773'
$ws.Range('E12').Value = '  -13.41%  '
$ws.Range('D13').Value = '14.75'
$ws.Range('E13').Value = '  -22.92%  '
$ws.Range('D14').Value = '4.452'
$ws.Range('E14').Value = '  -12.04%  '
$ws.Range('D15').Value = '77.19'
$ws.Range('E15').Value = '  -14.01%  '
$ws.Range('D16').Value = '0.5658'
$ws.Range('E16').Value = '  -27.59%  '
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').Value = '25.731.69'
$ws.Range('E19').Value = '  -3.79%  '
$ws.Range('D20').Value = '11.45'
$ws.Range('E20').Value = '  -19.24%  '
$ws.Range('D21').Value = '0.000006560'
$ws.Range('E21').Value = '  -18.36%  '
$ws.Range('D22').Value = '1.963.78'
$ws.Range('D23').Value = '3.960'
$ws.Range('E23').Value = '  -14.98%  '
$ws.Range('D24').Value = '4.997'
$ws.Range('E24').Value = '  -17.66%  '
$ws.Range('D25').Value = '7.847'
$ws.Range('E25').Value = '  -16.99%  '
$ws.Range('D26').Value = '135.71'
$ws.Range('E26').Value = '  -4.99%  '
$ws.Range('D27').Value = '1.476'
$ws.Range('E27').Value = '  -13.14%  '
$ws.Range('D28').Value = '1.814'
$ws.Range('E28').Value = '  -18.72%  '
$ws.Range('E29').Value = '  -14.68%  '
$ws.Range('D30').Value = '102.01'
$ws.Range('E30').Value = '  -8.76%  '
$ws.Range('D31').Value = '3.758'
$ws.Range('E31').Value = '  -12.67%  '
$ws.Range('D32').Value = '0.08054'
$ws.Range('E32').Value = '  -8.00%  '
$ws.Range('D33').Value = '3.311'
$ws.Range('E33').Value = '  -19.88%  '
$ws.Range('D34').Value = '0.04387'
$ws.Range('E34').Value = '  -10.28%  '
$ws.Range('D35').Value = '0.9999'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = '2.607'
$ws.Range('E36').Value = '  -9.12%  '
$ws.Range('D37').Value = '0.9731'
$ws.Range('E37').Value = '  -15.41%  '
$ws.Range('D38').Value = '0.5990'
$ws.Range('E38').Value = '  -19.08%  '
$ws.Range('E39').Value = '  -14.07%  '
$ws.Range('E40').Value = '  -14.80%  '
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').Value = '103.43'
$ws.Range('E42').Value = '  -5.59%  '
$ws.Range('D43').Value = '0.01490'
$ws.Range('E43').Value = '  -15.37%  '
$ws.Range('D44').Value = '0.7552'
$ws.Range('E44').Value = '  -17.48%  '
$ws.Range('D45').Value = '5.148'
$ws.Range('E45').Value = '  -13.33%  '
$ws.Range('D46').Value = '0.3703'
$ws.Range('E46').Value = '  -23.86%  '
$ws.Range('E47').Value = '  -12.58%  '
$ws.Range('D48').Value = '0.1064'
$ws.Range('E48').Value = '  -15.52%  '
$ws.Range('D49').Value = '30.06'
$ws.Range('E49').Value = '  -14.42%  '
$ws.Range('D50').Value = '52.42'
$ws.Range('E50').Value = '  -13.47%  '
$ws.Range('D51').Value = '5.859'
$ws.Range('E51').Value = '  -24.41%  '

# Restore the default cell style so formatting matches the original sheet.
$priceVolRange.Style = 'Normal'
